$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit performs a cyclic rotation of the data "records" (columns A, B, D, E,
# F, G, H, M, Q, R) that live in rows 3, 5, 6 and 7, while leaving every other
# column (which is identical across these rows) untouched:
#   new row 3 <- old row 6
#   new row 5 <- old row 3
#   new row 6 <- old row 7
#   new row 7 <- old row 5
#
# Read every source value first (via Value2, since the Value getter is not
# reliable in this runtime) so that later writes don't clobber values we
# still need to read.

$cols = @("A","B","D","E","F","G","H","Q","R")

$row3 = @{}
$row5 = @{}
$row6 = @{}
$row7 = @{}

foreach ($c in $cols) {
    $row3[$c] = $ws.Range("${c}3").Value2
    $row5[$c] = $ws.Range("${c}5").Value2
    $row6[$c] = $ws.Range("${c}6").Value2
    $row7[$c] = $ws.Range("${c}7").Value2
}

$m6 = $ws.Range("M6").Value2
$m7 = $ws.Range("M7").Value2

# Apply the rotation: new3 = old6, new5 = old3, new6 = old7, new7 = old5
foreach ($c in $cols) {
    $ws.Range("${c}3").Value2 = $row6[$c]
    $ws.Range("${c}5").Value2 = $row3[$c]
    $ws.Range("${c}6").Value2 = $row7[$c]
    $ws.Range("${c}7").Value2 = $row5[$c]
}

# Column M: row 6 had nothing, row 7 had "färska gnagspår". After the
# rotation row 6 receives row 7's old "M" value and row 7 ends up empty.
$ws.Range("M6").Value2 = $m7
$ws.Range("M7").Value2 = $m6
